# ----------------------------------------------------------------------------
# Rebuilds "my_comptability_sheet.xlsx" per the target commit:
#  - extends the income/expense table (Tableau1) with Epargne / Investissement
#    / "Depenses Fixe" / Reste (calculated) columns
#  - relocates the Actif/Passif table (Tableau3) from F:G to I:J
#  - adds a "Template" sheet that is a duplicate of the (now restructured)
#    "2023" sheet, complete with its own two tables
#  - touches up merges / headers / column widths / selection to match
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlPasteFormats = -4122

function Build-Sheet($ws) {
    # ---- months column (already present, but make sure) ------------------
    $months = @("Janvier","Février","Mars","Avril","Mai","Juin","Juillet","Août","Septembre","Octobre","Novembre","Décembre")
    for ($i = 0; $i -lt $months.Length; $i++) {
        $ws.Cells.Item(5 + $i, 2).Value = $months[$i]
    }

    # ---- tables: relocate Actif/Passif out of F:G, then grow Tableau1 ------
    $t1 = $null
    $t2 = $null
    foreach ($t in $ws.ListObjects) {
        if ($t.Range.Column -eq 2) { $t1 = $t }
        if ($t.Range.Column -eq 6 -or $t.Range.Column -eq 9) { $t2 = $t }
    }

    # move Actif/Passif table out of the way first (it currently sits on F:G)
    $ws.Range("I4").Value = "Actif"
    $ws.Range("J4").Value = "Passif"
    $t2.Resize($ws.Range("I4:J16"))

    # now F:G is free of the old table -> clear the stale header cells
    $ws.Range("F4:G16").ClearContents()

    # grow Tableau1 across the freed columns
    $t1.Resize($ws.Range("B4:G16"))

    # ---- table headers on row 4 -------------------------------------------
    $ws.Range("B4").Value = "Mois"
    $ws.Range("C4").Value = "Revenus"
    $ws.Range("D4").Value = "Epargne"
    $ws.Range("E4").Value = "Investissement"
    $ws.Range("F4").Value = "Depenses Fixe"
    $ws.Range("G4").Value = "Reste"
    $ws.Range("I4").Value = "Actif"
    $ws.Range("J4").Value = "Passif"

    # ---- "Reste" calculated column (shared formula) -----------------------
    $ws.Range("G5:G16").Formula = "= C5 - D5 -E5 - F5"

    # ---- headers on row 2 (band titles) ------------------------------------
    $ws.Range("B2").Value = "état des résultats"
    $ws.Range("I2").Value = "Bilan"

    # ---- totals row ----------------------------------------------------------
    $ws.Range("B18").Value = "Sommes Annuel"
    $ws.Range("C18").Formula = "=SUM(C5:C16)"
    $ws.Range("D18").Formula = "=SUM(D5:D16)"

    $ws.Range("C18").Copy()
    $ws.Range("E18").PasteSpecial($xlPasteFormats)
    $ws.Range("D18").Copy()
    $ws.Range("F18").PasteSpecial($xlPasteFormats)
    $ws.Range("D18").Copy()
    $ws.Range("G18").PasteSpecial($xlPasteFormats)
    $ws.Range("C18").Copy()
    $ws.Range("I18").PasteSpecial($xlPasteFormats)
    $ws.Range("D18").Copy()
    $ws.Range("J18").PasteSpecial($xlPasteFormats)

    $ws.Range("E18").Formula = "=SUM(E5:E16)"
    $ws.Range("F18").Formula = "=SUM(F5:F16)"
    $ws.Range("G18").Formula = "=SUM(G5:G16)"
    $ws.Range("I18").Formula = "=SUM(I5:I16)"
    $ws.Range("J18").Formula = "=SUM(J5:J16)"

    # ---- header band alignment / merges ------------------------------------
    $ws.Range("E2:G2").HorizontalAlignment = $xlCenter
    $ws.Range("J2").HorizontalAlignment = $xlCenter
    $ws.Range("B2:G2").Merge()
    $ws.Range("I2:J2").Merge()

    # ---- column widths for the two new columns -----------------------------
    $ws.Columns.Item(5).ColumnWidth = 16.02
    $ws.Columns.Item(6).ColumnWidth = 15.31

    return $t1.Name, $t2.Name
}

$ws1 = $wb.ActiveSheet
Build-Sheet $ws1

# ---- duplicate the sheet as "Template" -----------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Template"
Build-Sheet $ws2

$t1b = $ws2.ListObjects.Item(1)
$t1b.Name = "Tableau13"
$t2b = $ws2.ListObjects.Item(2)
$t2b.Name = "Tableau35"

# ---- selection / window bits ----------------------------------------------
$ws1.Range("A20").Select()

$wb.Save()
